# Sudoku "TestVector" sheet: fill in the missing digits in red.
$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> TestVector --------------------------------------------
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "TestVector"

# --- Fill in the missing Sudoku digits (shown in red) -----------------------
# Map of cell -> value for every cell that was previously blank.
$values = [ordered]@{
    "B2" = 7;  "D2" = 9;  "F2" = 2;  "H2" = 3;  "J2" = 6;
    "B3" = 4;  "D3" = 6;  "G3" = 8;  "I3" = 7;
    "B4" = 3;  "C4" = 1;  "E4" = 5;  "G4" = 6;  "H4" = 8;
    "E5" = 7;  "F5" = 4;  "H5" = 5;  "I5" = 6;  "J5" = 8;
    "B6" = 8;  "C6" = 6;  "E6" = 9;  "G6" = 5;  "J6" = 3;
    "C7" = 4;  "D7" = 3;  "F7" = 6;  "H7" = 7;  "I7" = 1;  "J7" = 9;
    "B8" = 6;  "D8" = 8;  "F8" = 3;  "H8" = 9;  "I8" = 2;
    "B9" = 1;  "F9" = 5;  "I9" = 8;  "J9" = 7;
    "D10" = 5; "E10" = 6; "G10" = 7; "H10" = 4; "I10" = 3
}

foreach ($ref in $values.Keys) {
    $cell = $ws.Range($ref)
    $cell.Value = $values[$ref]
    $cell.Font.Name = "Courier New"
    $cell.Font.Size = 14
    $cell.Font.Bold = $true
    $cell.Font.Color = 255
}

# --- Update the saved selection ---------------------------------------------
$ws.Range("I11").Select() | Out-Null
